$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Big Meadow Creek" data row (original row 3); rows below shift up.
$ws.Rows.Item(3).Delete()

# Remove the two trailing duplicate rows (originally rows 12-13, now rows 11-12
# after the shift above) that duplicated the Salmon/Aeneas entries.
$ws.Range("A11:F12").EntireRow.Delete()
